$d = $word.ActiveDocument

# Locate and remove the "Beach Day Studios - Game Programmer ... May 2020 - Present"
# paragraph (the newest job entry) that was added to the Work Experience section,
# while leaving the "Work Experience" heading and the following
# "Beach Day Studios - Game Programmer Internship" paragraph untouched.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    if ($text -like "Beach Day Studios*Game Programmer*" -and $text -like "*May 2020 - Present*") {
        $p.Range.Delete()
    }
}
